$wb = $excel.ActiveWorkbook

# --- Sheet1: raw data sheet ---
$ws1 = $wb.Worksheets.Item("Sheet1")

# The "id" column (A) is no longer populated for the data rows; remove those
# cells entirely (contents + formatting) so they disappear from the sheet XML.
$ws1.Range("A2:A14").Clear()

# The affinity_passive_id column (C) for Cold / Poison / Blood now stores a
# descriptive text string instead of a numeric id.
$ws1.Range("C11").Value = "Causes frost buildup ()"
$ws1.Range("C12").Value = "Causes poison buildup ()"
$ws1.Range("C13").Value = "Causes blood loss buildup ()"

# --- "table" sheet: presentation/view sheet ---
$ws2 = $wb.Worksheets.Item("table")

$ws2.Range("C11").Value = "Causes frost buildup ()"
$ws2.Range("C12").Value = "Causes poison buildup ()"
$ws2.Range("C13").Value = "Causes blood loss buildup ()"
